# Update gh-pages to output generated at 456a3b4
# Applies updated scrape counts (column F) across the four sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibit.Range("F3").Value  = 4707
$wsExhibit.Range("F4").Value  = 612
$wsExhibit.Range("F6").Value  = 1880
$wsExhibit.Range("F7").Value  = 556
$wsExhibit.Range("F12").Value = 1148
$wsExhibit.Range("F20").Value = 32
$wsExhibit.Range("F21").Value = 32
$wsExhibit.Range("F33").Value = 4313

# 演出 (sheet2)
$wsShow.Range("F21").Value = 2
$wsShow.Range("F32").Value = 5
$wsShow.Range("F39").Value = 40

# 本地生活 (sheet3)
$wsLocal.Range("F4").Value = 1354

# 全部类型 (sheet4)
$wsAll.Range("F3").Value  = 1354
$wsAll.Range("F9").Value  = 4707
$wsAll.Range("F10").Value = 612
$wsAll.Range("F12").Value = 1880
$wsAll.Range("F13").Value = 556
$wsAll.Range("F19").Value = 1148
$wsAll.Range("F30").Value = 32
$wsAll.Range("F47").Value = 5
$wsAll.Range("F49").Value = 4313
